# [ADDITIONAL SCRAPING] add player info + extra batting data, rework match-link
# columns into plain MATCH_CODE columns across the existing sheets.
#
# NOTE: worksheet handles returned by Worksheets.Item(...) are positional
# proxies, not stable object references - after any sheet is added/moved the
# same variable can end up pointing at a different sheet. So every sheet is
# re-fetched *by name* immediately before it is used.

$wb = $excel.ActiveWorkbook

function Style-Header($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 1) Create the new sheets first (structure), then fill them in afterwards.
# ---------------------------------------------------------------------------

# "Player Info" goes before the current first sheet ("ODI Batting").
$beforeTarget = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($beforeTarget)
$playerInfo.Name = "Player Info"

# "ODI Batting Extra" goes after the current last sheet ("ODI Bowling").
$afterTarget = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $afterTarget)
$extra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# 2) "Player Info" contents.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $piHeaders[$c - 1]
}
Style-Header $playerInfo.Range("A1:D1")

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3925"
$playerInfo.Range("B2").Value = "Jason Nazimuddin Mohammed"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ---------------------------------------------------------------------------
# 3) "ODI Batting" sheet - MATCH_CARD_LINK -> MATCH_CODE.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$battingCodes = @(
    "3360", "3853", "4001", "4004", "4005", "4017", "4018", "4019", "4040",
    "4043", "4046", "4051", "4052", "4053", "4056", "4057", "4068", "4070",
    "4072", "4073", "4075", "4100", "4101", "4102", "4144", "4148", "4179",
    "4180", "4443", "4445", "4447", "4449", "4450", "4451", "4483", "4484"
)

$odiBatting.Range("D1").Value = "MATCH_CODE"

$odiBatting.Range("D2:D37").NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $odiBatting.Cells.Item($i + 2, 4).Value = $battingCodes[$i]
}

# Rows where the player didn't bat (no INNING_NUMBER) lose their blank
# placeholder cell in column B entirely.
$blankInningRows = @(12, 13, 19, 26, 34)
foreach ($r in $blankInningRows) {
    $odiBatting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 4) "ODI Bowling" sheet - MATCH_CARD_LINK -> MATCH_CODE.
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

$bowlingCodes = @(
    "3853", "4001", "4004", "4017", "4019", "4100", "4144", "4179", "4443",
    "4445", "4447", "4449", "4450", "4451", "4483", "4484"
)

$odiBowling.Range("B1").Value = "MATCH_CODE"

$odiBowling.Range("B2:B17").NumberFormat = "@"
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $odiBowling.Cells.Item($i + 2, 2).Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------------
# 5) "ODI Batting Extra" contents.
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $extra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
Style-Header $extra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4068", 5,    "0", "1", "8.82%",  "NO"),
    @("4070", $null, $null, $null, $null, "NO"),
    @("4072", 5,    "3", "1", "15.51%", "NO"),
    @("4073", 5,    "4", "0", "12.92%", "NO"),
    @("4075", 5,    "1", "1", "8.68%",  "NO"),
    @("4100", 5,    "2", "0", "3.63%",  "NO"),
    @("4101", 5,    "2", "0", "14.88%", "NO"),
    @("4102", $null, $null, $null, $null, "NO"),
    @("4144", $null, $null, $null, $null, "NO"),
    @("4148", $null, $null, $null, $null, "NO"),
    @("4179", $null, $null, $null, $null, "NO"),
    @("4180", 5,    "1", "0", "4.43%",  "NO"),
    @("4443", $null, $null, $null, $null, "NO"),
    @("4445", 5,    "0", "0", "7.43%",  "NO"),
    @("4447", 5,    "1", "0", "9.60%",  "NO"),
    @("4449", 4,    "2", "0", "5.51%",  "NO"),
    @("4450", $null, $null, $null, $null, "NO"),
    @("4451", 3,    "1", "0", "2.90%",  "NO"),
    @("4483", 3,    "0", "0", "1.63%",  "NO"),
    @("4484", 4,    "2", "0", "5.76%",  "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = $i + 2
    $row = $extraRows[$i]

    # MATCH_CODE is always a literal code string, e.g. "4068".
    $codeCell = $extra.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    # BATTING_POSITION is a genuine number when known.
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }

    # NUM_4 / NUM_6 are literal digit strings, PERCENT_RUNS_OF_TOTAL keeps its
    # literal "%" suffix - only set these (as text) when the row has data.
    if ($null -ne $row[2]) {
        $c3 = $extra.Cells.Item($r, 3)
        $c3.NumberFormat = "@"
        $c3.Value = $row[2]

        $c4 = $extra.Cells.Item($r, 4)
        $c4.NumberFormat = "@"
        $c4.Value = $row[3]

        $c5 = $extra.Cells.Item($r, 5)
        $c5.NumberFormat = "@"
        $c5.Value = $row[4]
    }

    $extra.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Leave the active tab on the first sheet, as before the edit.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()
